$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(12).Insert()

$src1 = $ws.Range("D9:E9")
$src1.Copy()
$dst1 = $ws.Range("D12:E12")
$dst1.PasteSpecial(-4122)  # xlPasteFormats

$src2 = $ws.Range("F10:J10")
$src2.Copy()
$dst2 = $ws.Range("F12:J12")
$dst2.PasteSpecial(-4122)  # xlPasteFormats

$ws.Rows(12).EntireRow.AutoFit()
"done"
